$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '59.364.20'
Set-TextValue 'E2' '  +0.49%  '
Set-TextValue 'D3' '2.607.98'
Set-TextValue 'E3' '  +0.45%  '
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '546.50'
Set-TextValue 'E5' '  +4.74%  '
Set-TextValue 'E6' '  -0.11%  '
Set-TextValue 'E7' '  +0.28%  '
Set-TextValue 'D8' '0.567'
Set-TextValue 'E8' '  +0.05%  '
Set-TextValue 'D9' '6.47'
Set-TextValue 'E9' '  -1.08%  '
Set-TextValue 'E10' '  +1.72%  '
Set-TextValue 'E11' '  +0.82%  '
Set-TextValue 'D12' '0.135'
Set-TextValue 'E12' '  +1.49%  '
Set-TextValue 'D13' '3.068.97'
Set-TextValue 'E13' '  +0.32%  '
Set-TextValue 'D14' '59.290.67'
Set-TextValue 'E14' '  +0.34%  '
Set-TextValue 'D15' '20.60'
Set-TextValue 'E15' '  +0.07%  '
Set-TextValue 'B16' 'ShibaInu'
Set-TextValue 'C16' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D16' '0.0000134'
Set-TextValue 'E16' '  +0.48%  '
Set-TextValue 'B17' 'WrappedEther'
Set-TextValue 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '2.597.34'
Set-TextValue 'E17' '  -1.12%  '
Set-TextValue 'D18' '343.39'
Set-TextValue 'E18' '  +1.03%  '
Set-TextValue 'E19' '  +0.48%  '
Set-TextValue 'D20' '10.13'
Set-TextValue 'E20' '  -1.02%  '
Set-TextValue 'E21' '  -2.47%  '
Set-TextValue 'E22' '  +0.28%  '
Set-TextValue 'D23' '67.46'
Set-TextValue 'E23' '  +1.73%  '
Set-TextValue 'E24' '  -0.87%  '
Set-TextValue 'E25' '  +0.43%  '
Set-TextValue 'E26' '  +0.14%  '
Set-TextValue 'D27' '7.22'
Set-TextValue 'E27' '  +1.12%  '
Set-TextValue 'B28' 'USDe'
Set-TextValue 'C28' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  +0.14%  '
Set-TextValue 'B29' 'PEPE'
Set-TextValue 'C29' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D29' '0.0₃0738'
Set-TextValue 'E29' '  +1.23%  '
Set-TextValue 'D30' '1.71'
Set-TextValue 'E30' '  +8.88%  '
Set-TextValue 'D31' '5.82'
Set-TextValue 'E31' '  -2.66%  '
Set-TextValue 'D32' '18.78'
Set-TextValue 'E32' '  -0.12%  '
Set-TextValue 'D33' '149.29'
Set-TextValue 'E33' '  +0.19%  '
Set-TextValue 'E34' '  -0.77%  '
Set-TextValue 'D35' '37.10'
Set-TextValue 'E35' '  +2.09%  '
Set-TextValue 'E36' '  -1.85%  '
Set-TextValue 'E37' '  -0.22%  '
Set-TextValue 'E38' '  +0.01%  '
Set-TextValue 'E39' '  -2.30%  '
Set-TextValue 'D40' '3.55'
Set-TextValue 'E40' '  -0.26%  '
Set-TextValue 'B41' 'Bittensor'
Set-TextValue 'C41' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D41' '277.94'
Set-TextValue 'E41' '  +0.13%  '
Set-TextValue 'B42' 'FirstDigitalUSD'
Set-TextValue 'C42' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D42' '0.999'
Set-TextValue 'E42' '  +0.44%  '
Set-TextValue 'D43' '0.598'
Set-TextValue 'E43' '  +1.08%  '
Set-TextValue 'E44' '  +0.12%  '
Set-TextValue 'D45' '0.0957'
Set-TextValue 'E45' '  +0.01%  '
Set-TextValue 'E46' '  +0.38%  '
Set-TextValue 'D47' '1.948.56'
Set-TextValue 'E47' '  -1.97%  '
Set-TextValue 'E48' '  +0.83%  '
Set-TextValue 'B49' 'InjectiveProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D49' '18.35'
Set-TextValue 'E49' '  +1.15%  '
Set-TextValue 'B50' 'RenderToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D50' '4.52'
Set-TextValue 'E50' '  -2.35%  '
Set-TextValue 'D51' '111.06'
Set-TextValue 'E51' '  -2.68%  '
